$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the Price (D) and Volume(1h) (E) columns with the latest scraped
# market data. Price-column values are set to Text format first so that
# numeric-looking strings (e.g. "1.00", "604.89") are stored verbatim
# instead of being normalized by automatic number detection.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.372.12'
$ws.Range('E2').Value = '  +0.22%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.551.81'
$ws.Range('E3').Value = '  +0.45%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '604.89'
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.44'
$ws.Range('E6').Value = '  +0.49%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.550.12'
$ws.Range('E7').Value = '  +0.46%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.493'
$ws.Range('E9').Value = '  +2.72%  '
$ws.Range('E10').Value = '  -0.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.90'
$ws.Range('E11').Value = '  -2.00%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.413'
$ws.Range('E12').Value = '  +0.40%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.155.27'
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('E14').Value = '  +0.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '29.93'
$ws.Range('E15').Value = '  -1.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.554.66'
$ws.Range('E16').Value = '  +0.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.461.93'
$ws.Range('E17').Value = '  +0.22%  '
$ws.Range('E18').Value = '  +0.08%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.51'
$ws.Range('E19').Value = '  +5.54%  '
$ws.Range('E20').Value = '  -0.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.83'
$ws.Range('E21').Value = '  -0.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '430.96'
$ws.Range('E22').Value = '  +1.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.610'
$ws.Range('E23').Value = '  +1.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.60'
$ws.Range('E24').Value = '  +1.08%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.691.85'
$ws.Range('E25').Value = '  +0.49%  '
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('E27').Value = '  -0.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.51'
$ws.Range('E28').Value = '  +1.53%  '
$ws.Range('E29').Value = '  -0.33%  '
$ws.Range('E30').Value = '  -1.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.13%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '25.43'
$ws.Range('E32').Value = '  +0.48%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.546.61'
$ws.Range('E33').Value = '  +0.62%  '
$ws.Range('E34').Value = '  -2.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.153'
$ws.Range('E35').Value = '  -5.35%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.84'
$ws.Range('E36').Value = '  -0.07%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('E38').Value = '  -1.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.59'
$ws.Range('E39').Value = '  -0.63%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '175.61'
$ws.Range('E40').Value = '  +1.91%  '
$ws.Range('E41').Value = '  -1.28%  '
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('E43').Value = '  -0.48%  '
$ws.Range('E44').Value = '  +1.29%  '
$ws.Range('E45').Value = '  +1.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.999'
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.53'
$ws.Range('E47').Value = '  +5.46%  '
$ws.Range('E48').Value = '  -1.97%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.03'
$ws.Range('E49').Value = '  -3.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.14'
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '23.37'
$ws.Range('E51').Value = '  +3.37%  '
